# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# Concretely, for this workbook (sm_car_data_Aero_Coefficients.xlsx):
#  1. Update the "Sedan_Hamba" sheet's CD aero coefficient (H8) from 2.81 -> 1.98
#     and tighten its display format to 2 decimals (matches the other numeric rows).
#  2. Add a new "FSAE_Achilles" sheet (an FSAE car aero-coefficients entry), cloned
#     from the last existing sheet so it keeps identical layout/styles/tab color,
#     then overwritten with the FSAE-specific aero numbers.
#  3. Leave the other vehicle sheets' data untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sedan_Hamba: CD (row 8, column H) 2.81 -> 1.98, shown with 2 decimals.
# ---------------------------------------------------------------------------
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Range("H8").Value = 1.98
$sedanHamba.Range("H8").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 2) New sheet "FSAE_Achilles", appended after "Trailer_Kumanzi" (the last
#    sheet). Copying an existing sheet keeps its column widths, tab color and
#    per-cell styles identical, so only the data cells need to be rewritten.
# ---------------------------------------------------------------------------
$sourceCount = $wb.Worksheets.Count
$sourceSheet = $wb.Worksheets.Item($sourceCount)
$sourceSheet.Copy($null, $sourceSheet)

$fsae = $wb.Worksheets.Item($wb.Worksheets.Count)
$fsae.Name = "FSAE_Achilles"

$fsae.Range("H3").Value = "FSAE_Achilles"
$fsae.Range("H5").Value = -2.5
$fsae.Range("H6").Value = 1
$fsae.Range("H7").Value = 1.225
$fsae.Range("H8").Value = 1.2
$fsae.Range("F9").Value = -0.8
$fsae.Range("G9").Value = 0
$fsae.Range("H9").Value = 0.6

# Restore the per-pane selection state on the untouched sheets / new sheet so
# the saved view matches what a human editor would have left behind.
$wb.Worksheets.Item("Sedan_HambaLG").Range("H8").Select()
$fsae.Range("G12").Select()

# Sedan_Hamba ends up the active tab after the edit.
$sedanHamba.Range("E18").Select()
